$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 209 (pushes existing rows 209:230 down to 210:231,
# carrying along cell formatting such as the date style on column D).
$ws.Rows(209).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(209, 1).Value = 10
$ws.Cells.Item(209, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(209, 3).Value = "La Araucanía"
$ws.Cells.Item(209, 4).Value = 45154
$ws.Cells.Item(209, 5).Value = 9
$ws.Cells.Item(209, 6).Value = 100114002
$ws.Cells.Item(209, 7).Value = "Camote"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 30
$ws.Cells.Item(209, 11).Value = 26000
$ws.Cells.Item(209, 12).Value = 26000
$ws.Cells.Item(209, 13).Value = 26000
$ws.Cells.Item(209, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(209, 15).Value = "Perú"
$ws.Cells.Item(209, 16).Value = 1444
$ws.Cells.Item(209, 17).Value = 18
$ws.Cells.Item(209, 18).Value = "Hortaliza"
